$p = $ppt.ActivePresentation

# 1) Update the "date" auto-field cached text on the Handout Master (en-US)
#    and Notes Master (es-ES) from 28/03/2021 to 05/04/2021.
$p.HandoutMaster.HeadersFooters.DateAndTime.Text = "05/04/2021"
$p.NotesMaster.HeadersFooters.DateAndTime.Text = "05/04/2021"

# 2) Remove the "Rectángulo 7" shape (the "Notebooks: ..." banner) from slide 2.
$s2 = $p.Slides.Item(2)
for ($i = $s2.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s2.Shapes.Item($i)
    if ($shp.Name -eq "Rectángulo 7") {
        $shp.Delete()
    }
}
